# Update view-interest counts (column F, 想去人数) across sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 135
$ws.Range("F5").Value = 840
$ws.Range("F7").Value = 623
$ws.Range("F8").Value = 1193
$ws.Range("F10").Value = 757
$ws.Range("F11").Value = 668
$ws.Range("F12").Value = 252
$ws.Range("F13").Value = 354
$ws.Range("F14").Value = 345
$ws.Range("F15").Value = 718
$ws.Range("F16").Value = 861
$ws.Range("F17").Value = 9527
$ws.Range("F18").Value = 568
$ws.Range("F19").Value = 568
$ws.Range("F21").Value = 309
$ws.Range("F24").Value = 243
$ws.Range("F25").Value = 1726
$ws.Range("F26").Value = 20
$ws.Range("F27").Value = 277
$ws.Range("F28").Value = 472
$ws.Range("F29").Value = 168
$ws.Range("F31").Value = 247
$ws.Range("F32").Value = 179
$ws.Range("F33").Value = 51
$ws.Range("F34").Value = 91
$ws.Range("F36").Value = 168
$ws.Range("F37").Value = 169
$ws.Range("F38").Value = 153
$ws.Range("F39").Value = 33
$ws.Range("F40").Value = 89

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F7").Value = 108
$ws.Range("F12").Value = 71
$ws.Range("F16").Value = 255

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 800

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 800
$ws.Range("F5").Value = 135
$ws.Range("F8").Value = 840
$ws.Range("F10").Value = 623
$ws.Range("F11").Value = 1193
$ws.Range("F14").Value = 108
$ws.Range("F15").Value = 757
$ws.Range("F16").Value = 668
$ws.Range("F17").Value = 252
$ws.Range("F18").Value = 345
$ws.Range("F20").Value = 861
$ws.Range("F21").Value = 9527
$ws.Range("F23").Value = 568
$ws.Range("F24").Value = 568
$ws.Range("F27").Value = 243
$ws.Range("F28").Value = 1726
$ws.Range("F29").Value = 20
$ws.Range("F30").Value = 277
$ws.Range("F31").Value = 472
$ws.Range("F32").Value = 168
$ws.Range("F33").Value = 71
$ws.Range("F34").Value = 71
$ws.Range("F38").Value = 247
$ws.Range("F39").Value = 179
$ws.Range("F40").Value = 51
$ws.Range("F41").Value = 91
$ws.Range("F44").Value = 168
$ws.Range("F47").Value = 169
$ws.Range("F48").Value = 153
